# Apply edit described by commit "case studies discussion and figures":
#  - Split the SMT/F-RELATIVE speedup row (row 14) formulas into two
#    shared-formula groups (B14:G14 and H14:P14) by re-entering them.
#  - Add a new "SMT/FINE" speedup row (row 15) with formulas B9/B5 etc,
#    using the default (unstyled) number format.
#  - Update the sheet view: zoom to 115% and move the selection to L14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Row 14 (SMT/F-RELATIVE): re-enter the formulas so they are grouped
# into two separate shared-formula ranges, matching the edited layout.
$ws.Range("B14:G14").Formula = "=B9/B6"
$ws.Range("H14:P14").Formula = "=H9/H6"

# --- New row 15 (SMT/FINE) ---
$ws.Range("A15").Value = "SMT/FINE"
$ws.Range("B15:P15").Formula = "=B9/B5"
# Row 14 uses the 0.000 numeric style (s="2"); the new row keeps the
# default/general "Normal" style, so explicitly reset it.
$ws.Range("B15:P15").Style = "Normal"

# --- Sheet view changes: zoom in to 115% and select L14 ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("L14").Select()
